$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.130.97'
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").Value = '1.876.16'
$ws.Range("E3").Value = '  -0.94%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.22'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.97%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4891'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.76%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2901'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06584'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.00%  '

$ws.Range("D10").Value = '1.876.64'
$ws.Range("E10").Value = '  -0.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.38'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07199'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6646'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.907'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '85.83'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("D16").Value = '30.078.95'
$ws.Range("E16").Value = '  -0.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007787'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.73'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.49%  '

$ws.Range("D20").Value = '2.120.50'
$ws.Range("E20").Value = '  -0.78%  '

$ws.Range("E21").Value = '  +0.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.748'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.53%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.829'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.176'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.86'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.97%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '143.26'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +7.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.94'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.877'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.49%  '

$ws.Range("E29").Value = '  +1.63%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.196'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08769'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.44%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.990'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.10%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05124'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7113'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.61%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.105'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.669'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.15%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01838'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +10.59%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.678'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.60%  '

$ws.Range("E39").Value = '  -4.48%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9227'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.04%  '

$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '103.88'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.62%  '

$ws.Range("E42").Value = '  +0.23%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.758'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.36%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4212'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.399'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1279'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05713'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.20%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '32.77'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.08%  '

$ws.Range("B49").Value = 'Decentraland'
$ws.Range("C49").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.3747'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.23%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.216'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '55.66'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.56%  '
